$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 244: brand new match row (Warta Poznan vs Korona Kielce) ---
$ws.Range("B244").Value = 6775592
$ws.Range("E244").Value = 45394.54166666666
$ws.Range("F244").Value = "Warta Poznan"
$ws.Range("G244").Value = "Korona Kielce"
$ws.Range("H244").Value = 1
$ws.Range("I244").Value = 0
$ws.Range("J244").Value = "H"
$ws.Range("K244").Value = 2.45
$ws.Range("L244").Value = 3.1
$ws.Range("M244").Value = 2.75
$ws.Range("N244").Value = 2.45
$ws.Range("O244").Value = 2.75
$ws.Range("P244").Value = 3.1
$ws.Range("Q244").Value = 0
$ws.Range("R244").Value = 1.75
$ws.Range("S244").Value = 2.125
$ws.Range("T244").Value = 1.75
$ws.Range("U244").Value = 1.875
$ws.Range("V244").Value = 1.975
$ws.Range("W244").Value = 1.45
$ws.Range("X244").Value = -1
$ws.Range("Y244").Value = -1
$ws.Range("Z244").Value = 0.75
$ws.Range("AA244").Value = -1
$ws.Range("AB244").Value = -1
$ws.Range("AC244").Value = 0.9750000000000001

# --- Row 245: brand new match row (Pogon Szczecin vs Ruch Chorzow) ---
$ws.Range("B245").Value = 6774471
$ws.Range("E245").Value = 45394.64583333334
$ws.Range("F245").Value = "Pogon Szczecin"
$ws.Range("G245").Value = "Ruch Chorzow"
$ws.Range("H245").Value = 5
$ws.Range("I245").Value = 0
$ws.Range("J245").Value = "H"
$ws.Range("K245").Value = 1.5
$ws.Range("L245").Value = 4.2
$ws.Range("M245").Value = 5.5
$ws.Range("N245").Value = 1.5
$ws.Range("O245").Value = 4.333
$ws.Range("P245").Value = 5.5
$ws.Range("Q245").Value = -1
$ws.Range("R245").Value = 1.825
$ws.Range("S245").Value = 2.025
$ws.Range("T245").Value = 3
$ws.Range("U245").Value = 1.95
$ws.Range("V245").Value = 1.9
$ws.Range("W245").Value = 0.5
$ws.Range("X245").Value = -1
$ws.Range("Y245").Value = -1
$ws.Range("Z245").Value = 0.825
$ws.Range("AA245").Value = -1
$ws.Range("AB245").Value = 0.95
$ws.Range("AC245").Value = -1

# --- Row 246: LKS Lodz vs Radomiak Radom (odds correction) ---
$ws.Range("B246").Value = 6775593
$ws.Range("E246").Value = 45396.3125
$ws.Range("F246").Value = "LKS Lodz"
$ws.Range("G246").Value = "Radomiak Radom"
$ws.Range("K246").Value = 3.5
$ws.Range("L246").Value = 3.4
$ws.Range("M246").Value = 2.05
$ws.Range("N246").Value = 3.6
$ws.Range("O246").Value = 3.4
$ws.Range("P246").Value = 2
$ws.Range("Q246").Value = 0.5
$ws.Range("R246").Value = 1.825
$ws.Range("S246").Value = 2.025
$ws.Range("T246").Value = 2.5
$ws.Range("U246").Value = 2
$ws.Range("V246").Value = 1.85
$ws.Range("W246").Value = 0
$ws.Range("X246").Value = 0
$ws.Range("Y246").Value = 0
$ws.Range("Z246").Value = 0
$ws.Range("AA246").Value = 0

# --- Row 247: Jagiellonia Bialystok vs Cracovia Krakow (odds correction) ---
$ws.Range("B247").Value = 6775589
$ws.Range("E247").Value = 45396.41666666666
$ws.Range("F247").Value = "Jagiellonia Bialystok"
$ws.Range("G247").Value = "Cracovia Krakow"
$ws.Range("K247").Value = 1.833
$ws.Range("L247").Value = 3.5
$ws.Range("M247").Value = 3.8
$ws.Range("N247").Value = 1.65
$ws.Range("O247").Value = 3.75
$ws.Range("P247").Value = 4.5
$ws.Range("Q247").Value = -0.75
$ws.Range("R247").Value = 1.875
$ws.Range("S247").Value = 1.975
$ws.Range("T247").Value = 2.75
$ws.Range("U247").Value = 2.025
$ws.Range("V247").Value = 1.825
$ws.Range("W247").Value = 0
$ws.Range("X247").Value = 0
$ws.Range("Y247").Value = 0
$ws.Range("Z247").Value = 0
$ws.Range("AA247").Value = 0

# --- Row 248: Gornik Zabrze vs Slask Wroclaw (odds correction) ---
$ws.Range("B248").Value = 6775588
$ws.Range("E248").Value = 45396.52083333334
$ws.Range("F248").Value = "Gornik Zabrze"
$ws.Range("G248").Value = "Slask Wroclaw"
$ws.Range("K248").Value = 2.25
$ws.Range("L248").Value = 3.2
$ws.Range("M248").Value = 3
$ws.Range("N248").Value = 2.25
$ws.Range("O248").Value = 3.2
$ws.Range("P248").Value = 3
$ws.Range("Q248").Value = -0.25
$ws.Range("R248").Value = 2.025
$ws.Range("S248").Value = 1.825
$ws.Range("T248").Value = 2.25
$ws.Range("U248").Value = 1.9
$ws.Range("V248").Value = 1.95
$ws.Range("W248").Value = 0
$ws.Range("X248").Value = 0
$ws.Range("Y248").Value = 0
$ws.Range("Z248").Value = 0
$ws.Range("AA248").Value = 0

# --- Row 249: Piast Gliwice vs Zaglebie Lubin (odds correction) ---
$ws.Range("B249").Value = 6775590
$ws.Range("E249").Value = 45397.58333333334
$ws.Range("F249").Value = "Piast Gliwice"
$ws.Range("G249").Value = "Zaglebie Lubin"
$ws.Range("K249").Value = 2.15
$ws.Range("L249").Value = 3.2
$ws.Range("M249").Value = 3.4
$ws.Range("N249").Value = 2.15
$ws.Range("O249").Value = 3.2
$ws.Range("P249").Value = 3.5
$ws.Range("Q249").Value = -0.25
$ws.Range("R249").Value = 1.825
$ws.Range("S249").Value = 2.025
$ws.Range("T249").Value = 2
$ws.Range("U249").Value = 1.8
$ws.Range("V249").Value = 2.05
$ws.Range("W249").Value = 0
$ws.Range("X249").Value = 0
$ws.Range("Y249").Value = 0
$ws.Range("Z249").Value = 0
$ws.Range("AA249").Value = 0

# --- Row 250: no longer present in the refreshed data, remove it ---
$ws.Rows(250).Delete()
